$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion summary text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$new = @"
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 2.05 = 7515.4 pesos
✅ 7515.4 pesos = 2.04 = 878.53 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
"@
$ws1.Range("A1").Value = $new.TrimEnd("`r", "`n")

# --- Sheet "tasas": update the rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 487
$ws2.Range("O10").Value = 3660
$ws2.Range("N12").Value = 3687
$ws2.Range("O12").Value = 431
